$d = $word.ActiveDocument

# Locate the final (empty) paragraph in the document body - the new content
# is inserted immediately before it, matching the target diff.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionPoint = $lastPara.Range
$insertionPoint.Collapse(1)  # wdCollapseStart

$xmlFragment = @'
    <w:p>
      <w:pPr>
        <w:spacing w:before="60" w:after="100" w:afterAutospacing="1"/>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t xml:space="preserve">b1, b2, b3, b4 = </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t>Bike.new</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t>Bike.new</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t>Bike.new</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t>Bike.new</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:before="60" w:after="100" w:afterAutospacing="1"/>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t>b</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t>1.report</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t>_broken</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:before="60" w:after="100" w:afterAutospacing="1"/>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t>b</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t>2.report</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t>_broken</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:before="60" w:after="100" w:afterAutospacing="1"/>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t xml:space="preserve">s = </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t>DockingStation.new</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:before="60" w:after="100" w:afterAutospacing="1"/>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t>v</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t xml:space="preserve"> = </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t>Van.new</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:before="60" w:after="100" w:afterAutospacing="1"/>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t>s.dock</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t>(b1)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t xml:space="preserve"> &gt;&gt;</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:before="60" w:after="100" w:afterAutospacing="1"/>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t>s.send</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t>_broken_to</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
        <w:t>(v)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:before="60" w:after="100" w:afterAutospacing="1"/>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:color w:val="000000" w:themeColor="text1"/>
        </w:rPr>
      </w:pPr>
    </w:p>
'@

$flatOpc = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$xmlFragment</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

[void]$insertionPoint.InsertXML($flatOpc)
